$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 687
$ws.Range("F4").Value = 1994
$ws.Range("F5").Value = 5834
$ws.Range("F6").Value = 1644
$ws.Range("F7").Value = 174
$ws.Range("F8").Value = 3294
$ws.Range("F11").Value = 1383
$ws.Range("F12").Value = 4592
$ws.Range("F13").Value = 1096
$ws.Range("F14").Value = 1729
$ws.Range("F18").Value = 57
$ws.Range("F24").Value = 20
$ws.Range("F27").Value = 216
$ws.Range("F31").Value = 96
$ws.Range("F32").Value = 211
$ws.Range("F33").Value = 413
$ws.Range("F36").Value = 1763
$ws.Range("F37").Value = 2267
$ws.Range("F38").Value = 1059
$ws.Range("F42").Value = 389
$ws.Range("F43").Value = 42
$ws.Range("F44").Value = 680
$ws.Range("F47").Value = 426

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 22
$ws.Range("F11").Value = 160

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 687
$ws.Range("F5").Value = 1994
$ws.Range("F6").Value = 5834
$ws.Range("F7").Value = 1644
$ws.Range("F8").Value = 174
$ws.Range("F9").Value = 3294
$ws.Range("F11").Value = 1383
$ws.Range("F12").Value = 4592
$ws.Range("F13").Value = 1729
$ws.Range("F19").Value = 57
$ws.Range("F22").Value = 160
$ws.Range("F26").Value = 216
$ws.Range("F30").Value = 96
$ws.Range("F31").Value = 211
$ws.Range("F33").Value = 1763
$ws.Range("F34").Value = 2267
$ws.Range("F35").Value = 1059
$ws.Range("F41").Value = 389
$ws.Range("F42").Value = 680
$ws.Range("F44").Value = 426
